$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new header row at the top, shifting the existing keyword rows down
$ws.Rows.Item(1).Insert()

# Set the new header cell value and make it bold
$ws.Range("A1").Value = "Schlüsselwörter"
$ws.Range("A1").Font.Bold = $true

# Update the page setup (paper size / orientation), as Excel recorded on save
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Restore the selection to where the author last left it
$ws.Range("A6").Select() | Out-Null
